$d = $word.ActiveDocument

# The second paragraph (right after the title) is currently empty.
# Fill it with the centered, bold, red, underlined warning text.
$p = $d.Paragraphs.Item(2)
$r = $p.Range

# Seed the paragraph-mark's complex-script-bold flag before the text is
# (re)written, so the new run inherits it too once we re-apply below.
$r.Font.BoldBi = $true

$r.Text = "Make sure Visual Studio mode is set to x64 and not Win32"

$p.Alignment = 1  # wdAlignParagraphCenter

$r2 = $p.Range
$r2.Font.Bold = $true
$r2.Font.BoldBi = $true
$r2.Font.Color = 255        # wdColorRed (0xFF0000 -> 255)
$r2.Font.Size = 12
$r2.Font.SizeBi = 12
$r2.Font.Underline = 1      # wdUnderlineSingle
